$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 11
$ws.Range("C3").Value = 3
$ws.Range("C5").Value = 13
$ws.Range("B6").Value = "<like>"
$ws.Range("C6").Value = 8
$ws.Range("C7").Value = 11
$ws.Range("C8").Value = 8
$ws.Range("C9").Value = 11
$ws.Range("C10").Value = 13
$ws.Range("C11").Value = 9
$ws.Range("C12").Value = 8
$ws.Range("C13").Value = 10
$ws.Range("C14").Value = 13
$ws.Range("C16").Value = 8
$ws.Range("C18").Value = 6
